# Update the "Forecast Comparison" sheet:
#  - Insert a new "Week_Start_Date" column between "Week" and "ASIN"
#  - Normalize the Week labels (W01 -> W1, ... W09 -> W9; W10+ unchanged)
#  - Populate the new Week_Start_Date column with each week's start date (as text)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column B, shifting ASIN..is_holiday_week one column to the right.
$ws.Columns("B:B").Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Week labels (column A) and the corresponding week-start dates (new column B).
$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
$startDates = @("2025-01-05","2025-01-12","2025-01-19","2025-01-26","2025-02-02","2025-02-09","2025-02-16","2025-02-23","2025-03-02","2025-03-09","2025-03-16","2025-03-23","2025-03-30","2025-04-06","2025-04-13","2025-04-20")

for ($i = 0; $i -lt $weeks.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $weeks[$i]
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $startDates[$i]
}
